$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Blouse,Jumpsuit",
    "Kaftan,Tee",
    "Tee,Parka",
    "Kaftan,Tee",
    "Kaftan,Tee",
    "Dress,Jumpsuit",
    "Jumpsuit,Blouse",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Blazer",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Halter",
    "Jumpsuit,Dress",
    "Jumpsuit,Halter",
    "Jumpsuit,Halter",
    "Jumpsuit,Halter",
    "Jumpsuit,Halter",
    "Jumpsuit,Blouse",
    "Jumpsuit,Halter",
    "Jumpsuit,Dress",
    "Jumpsuit,Dress",
    "Tee,Parka",
    "Blazer,Jumpsuit",
    "Jumpsuit,Blouse",
    "Jumpsuit,Blouse",
    "Jumpsuit,Blouse",
    "Parka,Blouse",
    "Jumpsuit,Blouse",
    "Jumpsuit,Kaftan",
    "Halter,Blazer",
    "Jumpsuit,Kaftan"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
